$wb = $excel.ActiveWorkbook

# --- Rename MATCH_CARD_LINK -> MATCH_CODE and simplify link values to bare match codes ---

# ODI Batting sheet (column D)
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "4449"
$batting.Range("D3").Value = "4450"
$batting.Range("D4").Value = "4451"
$batting.Range("D5").Value = "4463"
$batting.Range("D6").Value = "4464"
$batting.Range("D7").Value = "4691"

# ODI Bowling sheet (column B)
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "4449"
$bowling.Range("B3").Value = "4691"

# --- Add new "Player Info" sheet as the first sheet ---

$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# match the bold / bordered / centered header style used on the other sheets
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$playerInfo.Range("A2").Value = "5950"
$playerInfo.Range("B2").Value = "Koralegedera Nadeeja Ashen Bandara"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

$playerInfo.Range("A1").Select() | Out-Null
